# Update the auto-updating "datetime6" date fields on the Handout Master
# and Notes Master from "novembre 22" to "avril 23" (the field recalculated
# to a later capture date).
$p = $ppt.ActivePresentation

$handoutMaster = $p.HandoutMaster
$handoutMaster.HeadersFooters.DateAndTime.Text = "avril 23"

$notesMaster = $p.NotesMaster
$notesMaster.HeadersFooters.DateAndTime.Text = "avril 23"

# Remove the two stray "PLATFORM" / "VIRTUALIZATION" label textboxes
# (Rectangle 162 / Rectangle 163) from the overview schematic slide.
$slide = $p.Slides.Item(1)

$platformShape = $slide.Shapes.Item("Rectangle 162")
$platformShape.Delete()

$virtualizationShape = $slide.Shapes.Item("Rectangle 163")
$virtualizationShape.Delete()
